$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.28617366666667
$ws.Range("H2").Value = 69.858521
$ws.Range("I2").Value = 0.2304887056246027
$ws.Range("J2").Value = 0.2304887056246027
$ws.Range("M2").Value = 23.28617366666667
$ws.Range("N2").Value = 69.858521
$ws.Range("O2").Value = 0.2304887056246027
$ws.Range("P2").Value = 0.2304887056246027
$ws.Range("Q2").Value = 542.2458840341601
$ws.Range("R2").Value = 4880.21295630744
$ws.Range("S2").Value = 0.05312504342050477
$ws.Range("T2").Value = 0.05312504342050475
$ws.Range("G3").Value = 23.28617366666667
$ws.Range("H3").Value = 69.858521
$ws.Range("I3").Value = 0.2304887056246027
$ws.Range("J3").Value = 0.2304887056246027
$ws.Range("O3").Value = 0.007098179626924059
$ws.Range("P3").Value = 0.007098179626924059
$ws.Range("Q3").Value = 16.699120576882
$ws.Range("R3").Value = 150.292085191938
$ws.Range("S3").Value = 0.001636050234500652
$ws.Range("T3").Value = 0.001636050234500652
$ws.Range("G4").Value = 23.28617366666667
$ws.Range("H4").Value = 69.858521
$ws.Range("I4").Value = 0.2304887056246027
$ws.Range("J4").Value = 0.2304887056246027
$ws.Range("O4").Value = 0.7624131147484733
$ws.Range("P4").Value = 0.7624131147484732
$ws.Range("Q4").Value = 1793.646991446746
$ws.Range("R4").Value = 16142.82292302072
$ws.Range("S4").Value = 0.1757276119695973
$ws.Range("T4").Value = 0.1757276119695973
$ws.Range("I5").Value = 0.007098179626924059
$ws.Range("J5").Value = 0.007098179626924059
$ws.Range("M5").Value = 23.28617366666667
$ws.Range("N5").Value = 69.858521
$ws.Range("O5").Value = 0.2304887056246027
$ws.Range("P5").Value = 0.2304887056246027
$ws.Range("Q5").Value = 16.699120576882
$ws.Range("R5").Value = 150.292085191938
$ws.Range("S5").Value = 0.001636050234500652
$ws.Range("T5").Value = 0.001636050234500652
$ws.Range("I6").Value = 0.007098179626924059
$ws.Range("J6").Value = 0.007098179626924059
$ws.Range("O6").Value = 0.007098179626924059
$ws.Range("P6").Value = 0.007098179626924059
$ws.Range("S6").Value = 0.00005038415401607977
$ws.Range("T6").Value = 0.00005038415401607977
$ws.Range("I7").Value = 0.007098179626924059
$ws.Range("J7").Value = 0.007098179626924059
$ws.Range("O7").Value = 0.7624131147484733
$ws.Range("P7").Value = 0.7624131147484732
$ws.Range("S7").Value = 0.005411745238407328
$ws.Range("T7").Value = 0.005411745238407327
$ws.Range("I8").Value = 0.7624131147484733
$ws.Range("J8").Value = 0.7624131147484732
$ws.Range("M8").Value = 23.28617366666667
$ws.Range("N8").Value = 69.858521
$ws.Range("O8").Value = 0.2304887056246027
$ws.Range("P8").Value = 0.2304887056246027
$ws.Range("Q8").Value = 1793.646991446746
$ws.Range("R8").Value = 16142.82292302072
$ws.Range("S8").Value = 0.1757276119695973
$ws.Range("T8").Value = 0.1757276119695973
$ws.Range("I9").Value = 0.7624131147484733
$ws.Range("J9").Value = 0.7624131147484732
$ws.Range("O9").Value = 0.007098179626924059
$ws.Range("P9").Value = 0.007098179626924059
$ws.Range("S9").Value = 0.005411745238407328
$ws.Range("T9").Value = 0.005411745238407327
$ws.Range("I10").Value = 0.7624131147484733
$ws.Range("J10").Value = 0.7624131147484732
$ws.Range("O10").Value = 0.7624131147484733
$ws.Range("P10").Value = 0.7624131147484732
$ws.Range("S10").Value = 0.5812737575404687
$ws.Range("T10").Value = 0.5812737575404685
